## Room Database - Raw.xlsx : update the "Joule" (B) column of sheet
## "Blad1" with the re-processed Nexus 5X measurements (rows 2-91),
## then leave the selection on D88 as it was when the file was last
## saved by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(111.14243999999999, 127.04219999999999, 128.13336000000001, 127.35396, 127.50984, 127.50984, 126.574559999999, 127.50984, 126.88632, 122.36579999999999, 127.8216, 127.19808, 126.73044, 127.04219999999999, 127.35396, 126.418679999999, 127.19808, 127.04219999999999, 127.04219999999999, 127.04219999999999, 127.19808, 127.97748, 127.04219999999999, 126.88632, 127.35396, 127.19808, 121.898159999999, 126.73044, 127.35396, 127.19808, 113.247828, 132.04414800000001, 115.91064, 116.067275999999, 133.92377999999999, 111.838104, 116.067275999999, 116.067275999999, 116.537184, 133.45387199999999, 116.380547999999, 115.440732, 133.297236, 133.297236, 133.76714399999901, 116.067275999999, 114.97082399999999, 115.28409600000001, 115.597368, 118.260179999999, 130.94769600000001, 131.73087599999999, 124.99552799999999, 132.35741999999999, 131.41760399999899, 131.26096799999999, 131.57424, 131.41760399999899, 131.88751199999999, 134.863596, 113.5611, 129.537972, 111.211559999999, 128.911428, 111.211559999999, 128.911428, 128.911428, 129.22469999999899, 129.06806399999999, 128.59815599999999, 129.85124399999901, 128.28488399999901, 129.06806399999999, 128.75479200000001, 129.06806399999999, 128.75479200000001, 128.75479200000001, 129.06806399999999, 129.06806399999999, 128.75479200000001, 128.75479200000001, 128.59815599999999, 124.99552799999999, 129.22469999999899, 128.911428, 128.59815599999999, 129.381336, 128.75479200000001, 128.911428, 128.59815599999999)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Restore the final selection recorded in the sheet view.
$ws.Range("D88").Select()
